$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 693.625
$ws.Range("I5").Value = 729.8
$ws.Range("J5").Value = 633.3333
$ws.Range("K5").Value = 729.8
$ws.Range("L5").Value = 633.3333
$ws.Range("M5").Value = -614.8
$ws.Range("N5").Value = -863.3333

$ws.Range("H58").Value = 621.44446
$ws.Range("I58").Value = 327.7143
$ws.Range("K58").Value = 983.1428999999999
$ws.Range("M58").Value = -833.1428999999999

$ws.Range("H97").Value = 1412.2
$ws.Range("J97").Value = 1412.2
$ws.Range("L97").Value = 4236.6
$ws.Range("N97").Value = -5228.6

$ws.Range("H115").Value = 3248.9092
$ws.Range("J115").Value = 20000
$ws.Range("L115").Value = 60000
$ws.Range("N115").Value = -63134

$ws.Range("H116").Value = 7866.1665
$ws.Range("I116").Value = 4878.8
$ws.Range("K116").Value = 4878.8
$ws.Range("M116").Value = -1436.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 25008268
$ws.Range("I37").Value = 25008268
$ws.Range("J37").Value = 0
$ws.Range("K37").Value = 25008268
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = -25007995
$ws.Range("N37").ClearContents()

$ws.Range("H45").Value = 24602.143
$ws.Range("I45").Value = 24824.889
$ws.Range("K45").Value = 24824.889
$ws.Range("M45").Value = -24447.889

$ws.Range("H55").Value = 15256

$ws.Range("H63").Value = 717.3333
$ws.Range("I63").Value = 717.3333
$ws.Range("K63").Value = 717.3333
$ws.Range("M63").Value = -31.33330000000001

$ws.Range("H66").Value = 717.3333
$ws.Range("I66").Value = 717.3333
$ws.Range("K66").Value = 3586.6665
$ws.Range("M66").Value = -154.6665000000003

$ws.Range("H74").Value = 2087.4707
$ws.Range("I74").Value = 1540.75
$ws.Range("K74").Value = 1540.75
$ws.Range("M74").Value = -666.75

$ws.Range("H77").Value = 2087.4707
$ws.Range("I77").Value = 1540.75
$ws.Range("K77").Value = 7703.75
$ws.Range("M77").Value = -3335.75

$ws.Range("H80").Value = 19998.334
$ws.Range("J80").Value = 19998.334
$ws.Range("L80").Value = 19998.334
$ws.Range("N80").Value = -21994.334

$ws.Range("H83").Value = 19998.334
$ws.Range("J83").Value = 19998.334
$ws.Range("L83").Value = 59995.00199999999
$ws.Range("N83").Value = -69979.00199999999

$ws.Range("H110").Value = 1942.8334
$ws.Range("I110").Value = 1942.8334
$ws.Range("K110").Value = 1942.8334
$ws.Range("M110").Value = 102.1666

$ws.Range("H122").Value = 14311.0625
$ws.Range("I122").Value = 18463.584
$ws.Range("K122").Value = 55390.75199999999
$ws.Range("M122").Value = -52940.75199999999

$ws.Range("H132").Value = 2247.7297
$ws.Range("I132").Value = 2316.7856
$ws.Range("K132").Value = 6950.3568
$ws.Range("M132").Value = -4420.3568

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H35").Value = 30000
$ws.Range("J35").Value = 30000
$ws.Range("L35").Value = 30000
$ws.Range("N35").Value = -30620

$ws.Range("H82").Value = 22179.572
$ws.Range("I82").Value = 15128.5
$ws.Range("J82").Value = 25000
$ws.Range("K82").Value = 15128.5
$ws.Range("L82").Value = 25000
$ws.Range("M82").Value = -14745.5
$ws.Range("N82").Value = -25766

$ws.Range("H85").Value = 22179.572
$ws.Range("I85").Value = 15128.5
$ws.Range("J85").Value = 25000
$ws.Range("K85").Value = 15128.5
$ws.Range("L85").Value = 25000
$ws.Range("M85").Value = -13802.5
$ws.Range("N85").Value = -27652

$ws.Range("H86").Value = 1761.4073
$ws.Range("I86").Value = 1648.174
$ws.Range("J86").Value = 2412.5
$ws.Range("K86").Value = 1648.174
$ws.Range("L86").Value = 2412.5
$ws.Range("M86").Value = -525.174
$ws.Range("N86").Value = -4658.5

$ws.Range("H89").Value = 1761.4073
$ws.Range("I89").Value = 1648.174
$ws.Range("J89").Value = 2412.5
$ws.Range("K89").Value = 8240.869999999999
$ws.Range("L89").Value = 12062.5
$ws.Range("M89").Value = -2624.869999999999
$ws.Range("N89").Value = -23294.5

$ws.Range("H107").Value = 1526.9166
$ws.Range("I107").Value = 874.8889
$ws.Range("J107").Value = 3483
$ws.Range("K107").Value = 874.8889
$ws.Range("L107").Value = 3483
$ws.Range("M107").Value = 1045.1111
$ws.Range("N107").Value = -7323

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 382
$ws.Range("J22").Value = 382
$ws.Range("L22").Value = 382
$ws.Range("N22").Value = -1082

$ws.Range("H35").Value = 1059.1428
$ws.Range("I35").Value = 402.33334
$ws.Range("J35").Value = 5000
$ws.Range("K35").Value = 402.33334
$ws.Range("L35").Value = 5000
$ws.Range("M35").Value = -108.33334
$ws.Range("N35").Value = -5588

$ws.Range("H105").Value = 2913.8276
$ws.Range("I105").Value = 2465.353
$ws.Range("K105").Value = 2465.353
$ws.Range("M105").Value = -718.3530000000001

$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("L106").Value = 0
$ws.Range("N106").ClearContents()

$ws.Range("H120").Value = 21500
$ws.Range("J120").Value = 21500
$ws.Range("L120").Value = 21500
$ws.Range("N120").Value = -28758

$ws.Range("H132").Value = 4833.905
$ws.Range("I132").Value = 4833.905
$ws.Range("K132").Value = 14501.715
$ws.Range("M132").Value = -11971.715

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3170.8125
$ws.Range("I129").Value = 1580
$ws.Range("J129").Value = 3537.923
$ws.Range("K129").Value = 4740
$ws.Range("L129").Value = 10613.769
$ws.Range("M129").Value = 260
$ws.Range("N129").Value = -20613.769

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 29999
$ws.Range("J57").Value = 29999
$ws.Range("L57").Value = 29999
$ws.Range("N57").Value = -31639

$ws.Range("H70").Value = 7109.8184
$ws.Range("I70").Value = 7055.4443
$ws.Range("K70").Value = 7055.4443
$ws.Range("M70").Value = -6785.4443

$ws.Range("H73").Value = 7109.8184
$ws.Range("I73").Value = 7055.4443
$ws.Range("K73").Value = 7055.4443
$ws.Range("M73").Value = -6119.4443

$ws.Range("H126").Value = 3154.375
$ws.Range("J126").Value = 3412
$ws.Range("L126").Value = 10236
$ws.Range("N126").Value = -15176

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("M4").ClearContents()

$ws.Range("H7").Value = 23299.867
$ws.Range("I7").Value = 35888.668
$ws.Range("K7").Value = 35888.668
$ws.Range("M7").Value = -35776.668

$ws.Range("H20").Value = 3430000
$ws.Range("I20").Value = 40000
$ws.Range("J20").Value = 5125000
$ws.Range("K20").Value = 40000
$ws.Range("L20").Value = 5125000
$ws.Range("M20").Value = -39774
$ws.Range("N20").Value = -5125452

$ws.Range("H22").Value = 1768.8
$ws.Range("I22").Value = 2046
$ws.Range("J22").Value = 1584
$ws.Range("K22").Value = 2046
$ws.Range("L22").Value = 1584
$ws.Range("M22").Value = -1751
$ws.Range("N22").Value = -2174

$ws.Range("H25").Value = 4950
$ws.Range("I25").Value = 4900
$ws.Range("K25").Value = 4900
$ws.Range("M25").Value = -4670

$ws.Range("H27").Value = 1768.8
$ws.Range("I27").Value = 2046
$ws.Range("J27").Value = 1584
$ws.Range("K27").Value = 2046
$ws.Range("L27").Value = 1584
$ws.Range("M27").Value = -1939
$ws.Range("N27").Value = -1798

$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 0
$ws.Range("K28").Value = 0
$ws.Range("M28").ClearContents()

$ws.Range("H37").Value = 0
$ws.Range("I37").Value = 0
$ws.Range("K37").Value = 0
$ws.Range("M37").ClearContents()

$ws.Range("H46").Value = 49446.777
$ws.Range("I46").Value = 72353.5
$ws.Range("J46").Value = 3633.3333
$ws.Range("K46").Value = 72353.5
$ws.Range("L46").Value = 3633.3333
$ws.Range("M46").Value = -72165.5
$ws.Range("N46").Value = -4009.3333

$ws.Range("H61").Value = 15168967
$ws.Range("I61").Value = 17546302
$ws.Range("K61").Value = 17546302
$ws.Range("M61").Value = -17546100

$ws.Range("H93").Value = 43340.25
$ws.Range("J93").Value = 113041
$ws.Range("L93").Value = 113041
$ws.Range("N93").Value = -115537

$ws.Range("H100").Value = 223888.33
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

$ws.Range("H113").Value = 15168967
$ws.Range("I113").Value = 17546302
$ws.Range("K113").Value = 17546302
$ws.Range("M113").Value = -17544132

$ws.Range("H122").Value = 8773.210999999999
$ws.Range("I122").Value = 10059.4
$ws.Range("J122").Value = 3950
$ws.Range("K122").Value = 30178.2
$ws.Range("L122").Value = 11850
$ws.Range("M122").Value = -27728.2
$ws.Range("N122").Value = -16750

$ws.Range("H126").Value = 23299.867
$ws.Range("I126").Value = 35888.668
$ws.Range("K126").Value = 107666.004
$ws.Range("M126").Value = -105196.004

$ws.Range("H136").Value = 4983
$ws.Range("I136").Value = 3724
$ws.Range("J136").Value = 7501
$ws.Range("K136").Value = 11172
$ws.Range("L136").Value = 22503
$ws.Range("M136").Value = -8622
$ws.Range("N136").Value = -27603

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 729207.7
$ws.Range("I14").Value = 784531.4
$ws.Range("J14").Value = 10000
$ws.Range("K14").Value = 784531.4
$ws.Range("L14").Value = 10000
$ws.Range("M14").Value = -784363.4
$ws.Range("N14").Value = -10336

$ws.Range("H133").Value = 80715
$ws.Range("J133").Value = 80715
$ws.Range("L133").Value = 80715
$ws.Range("N133").Value = -90835
